$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Map of (row, col) -> (old, new) text, matching the diff exactly.
# Using per-cell Find/Replace avoids ambiguity from duplicate values
# (e.g. "352×3=1056" appears in two different cells with different
# replacements).
$replacements = @(
    @{ Row = 1;  Col = 1; Old = "105×4=420";  New = "461×8=3688" },
    @{ Row = 1;  Col = 2; Old = "519×4=2076"; New = "403×9=3627" },
    @{ Row = 1;  Col = 3; Old = "352×3=1056"; New = "669×9=6021" },
    @{ Row = 1;  Col = 4; Old = "152×6=912";  New = "870×7=6090" },
    @{ Row = 1;  Col = 5; Old = "844×2=1688"; New = "606×5=3030" },

    @{ Row = 5;  Col = 1; Old = "615×3=1845"; New = "976×6=5856" },
    @{ Row = 5;  Col = 2; Old = "352×3=1056"; New = "463×4=1852" },
    @{ Row = 5;  Col = 3; Old = "720×8=5760"; New = "826×9=7434" },
    @{ Row = 5;  Col = 4; Old = "747×4=2988"; New = "730×7=5110" },
    @{ Row = 5;  Col = 5; Old = "611×4=2444"; New = "322×3=966"  },

    @{ Row = 10; Col = 1; Old = "567×3=1701"; New = "961×9=8649" },
    @{ Row = 10; Col = 2; Old = "818×8=6544"; New = "685×2=1370" },
    @{ Row = 10; Col = 3; Old = "232×2=464";  New = "759×7=5313" },
    @{ Row = 10; Col = 4; Old = "877×8=7016"; New = "647×2=1294" },
    @{ Row = 10; Col = 5; Old = "198×7=1386"; New = "311×6=1866" },

    @{ Row = 15; Col = 1; Old = "677×4=2708"; New = "925×7=6475" },
    @{ Row = 15; Col = 2; Old = "370×3=1110"; New = "545×4=2180" },
    @{ Row = 15; Col = 3; Old = "220×2=440";  New = "922×7=6454" },
    @{ Row = 15; Col = 4; Old = "493×2=986";  New = "617×9=5553" },
    @{ Row = 15; Col = 5; Old = "162×3=486";  New = "527×4=2108" },

    @{ Row = 20; Col = 1; Old = "402×5=2010"; New = "212×7=1484" },
    @{ Row = 20; Col = 2; Old = "610×8=4880"; New = "434×6=2604" },
    @{ Row = 20; Col = 3; Old = "755×8=6040"; New = "488×7=3416" },
    @{ Row = 20; Col = 4; Old = "434×8=3472"; New = "257×2=514"  },
    @{ Row = 20; Col = 5; Old = "374×5=1870"; New = "946×7=6622" }
)

foreach ($item in $replacements) {
    $cell = $t.Cell($item.Row, $item.Col)
    # Wrap = 0 (wdFindStop) and Replace = 1 (wdReplaceOne) keep the
    # substitution confined to this cell's Range even when the same
    # text appears elsewhere in the document (e.g. "352×3=1056" shows
    # up in two different cells with two different replacements).
    $cell.Range.Find.Execute($item.Old, $true, $false, $false, $false, $false,
                              $true, 0, $false, $item.New, 1)
}
